$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = '-587'
$ws.Range("A70").Style = "Normal"
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = '9/8/2025'
$ws.Range("B70").Style = "Normal"
$ws.Range("C70").Value = 'ARIAS 4048'
$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = '12'
$ws.Range("D70").Style = "Normal"
$ws.Range("E70").NumberFormat = "@"
$ws.Range("E70").Value = '809526164'
$ws.Range("E70").Style = "Normal"
$ws.Range("H70").Value = 'Cambiar 114 picada'
$ws.Range("K70").Value = 'Sin equipos'
$ws.Range("M70").Value = -58.488936
$ws.Range("N70").Value = -34.549005
# Row 71
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = '7229'
$ws.Range("A71").Style = "Normal"
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = '9/16/2025'
$ws.Range("B71").Style = "Normal"
$ws.Range("C71").Value = 'AZURDUY JUANA 2627'
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = '13'
$ws.Range("D71").Style = "Normal"
$ws.Range("E71").NumberFormat = "@"
$ws.Range("E71").Value = 'ICD30814490'
$ws.Range("E71").Style = "Normal"
$ws.Range("G71").Value = 'Pendiente de Traspaso PROPIO'
$ws.Range("H71").Value = 'Colocar columna para pedir traspaso de nodo propio'
$ws.Range("K71").Value = 'Nodo Teco'
$ws.Range("M71").Value = -58.469008
$ws.Range("N71").Value = -34.552083
# Row 72
$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = '2711'
$ws.Range("A72").Style = "Normal"
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = '9/22/2025'
$ws.Range("B72").Style = "Normal"
$ws.Range("C72").Value = 'RUIZ HUIDOBRO 3620'
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = '12'
$ws.Range("D72").Style = "Normal"
$ws.Range("E72").NumberFormat = "@"
$ws.Range("E72").Value = 'ICD30934235'
$ws.Range("E72").Style = "Normal"
$ws.Range("H72").Value = 'Trapaso de redes y desmonte'
$ws.Range("J72").Value = 'Desmonte'
$ws.Range("K72").Value = 'Sin equipos'
$ws.Range("M72").Value = -58.484082
$ws.Range("N72").Value = -34.549702
# Row 73
$ws.Range("A73").NumberFormat = "@"
$ws.Range("A73").Value = '7277'
$ws.Range("A73").Style = "Normal"
$ws.Range("C73").Value = 'LA FRONDA 1670'
$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value = '11'
$ws.Range("D73").Style = "Normal"
$ws.Range("E73").NumberFormat = "@"
$ws.Range("E73").Value = '809929900'
$ws.Range("E73").Style = "Normal"
$ws.Range("H73").Value = 'Columna inclinada'
$ws.Range("J73").Value = 'Cambio'
$ws.Range("M73").Value = -58.468984
$ws.Range("N73").Value = -34.61101
$ws.Range("O73").Value = 'Paternal'
# Row 74
$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = '-605'
$ws.Range("A74").Style = "Normal"
$ws.Range("C74").Value = 'Tronador 1836'
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = '12'
$ws.Range("D74").Style = "Normal"
$ws.Range("E74").NumberFormat = "@"
$ws.Range("E74").Value = '809972725'
$ws.Range("E74").Style = "Normal"
$ws.Range("H74").Value = 'Poste podrido'
$ws.Range("L74").Value = 'Poste'
$ws.Range("M74").Value = -58.470216
$ws.Range("N74").Value = -34.57369
$ws.Range("O74").Value = 'Colegiales'
# Row 75
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = '-606'
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value = '9/23/2025'
$ws.Range("B75").Style = "Normal"
$ws.Range("C75").Value = 'Av. Cordoba  5106'
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = '15'
$ws.Range("D75").Style = "Normal"
$ws.Range("E75").NumberFormat = "@"
$ws.Range("E75").Value = '809930303'
$ws.Range("E75").Style = "Normal"
$ws.Range("H75").Value = 'Por obra Pedro'
$ws.Range("I75").Value = 0
$ws.Range("L75").Value = 'Pasante'
$ws.Range("M75").Value = -58.435561
$ws.Range("N75").Value = -34.590765
$ws.Range("O75").Value = 'Palermo'
$ws.Range("P75").Value = 'Capital Sur'
# Row 76
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = '-607'
$ws.Range("A76").Style = "Normal"
$ws.Range("C76").Value = 'Thames 1195'
$ws.Range("E76").NumberFormat = "@"
$ws.Range("E76").Value = '809930323'
$ws.Range("E76").Style = "Normal"
$ws.Range("M76").Value = -58.435404
$ws.Range("N76").Value = -34.590934
# Row 77
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = '-608'
$ws.Range("A77").Style = "Normal"
$ws.Range("C77").Value = 'Av. Cordoba 5064'
$ws.Range("E77").NumberFormat = "@"
$ws.Range("E77").Value = '809930333'
$ws.Range("E77").Style = "Normal"
$ws.Range("M77").Value = -58.435062
$ws.Range("N77").Value = -34.5911
# Row 78
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = '-614'
$ws.Range("A78").Style = "Normal"
$ws.Range("B78").NumberFormat = "@"
$ws.Range("B78").Value = '9/25/2025'
$ws.Range("B78").Style = "Normal"
$ws.Range("C78").Value = 'O''Higgins 2471'
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = '13'
$ws.Range("D78").Style = "Normal"
$ws.Range("E78").NumberFormat = "@"
$ws.Range("E78").Value = '809972821'
$ws.Range("E78").Style = "Normal"
$ws.Range("H78").Value = 'Recambio'
$ws.Range("I78").Value = 1
$ws.Range("L78").Value = 'Poste'
$ws.Range("M78").Value = -58.455887
$ws.Range("N78").Value = -34.556394
$ws.Range("O78").Value = 'Saavedra'
$ws.Range("P78").Value = 'Capital Norte'
